{"js": "// Update the header date and every three-digit-by-one-digit multiplication\n// answer in the practice table, per the commit's regenerated data.\nconst replacements = [\n  [\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"],\n  [\"500\u00d75=2500\", \"338\u00d73=1014\"],\n  [\"586\u00d72=1172\", \"807\u00d76=4842\"],\n  [\"341\u00d74=1364\", \"574\u00d77=4018\"],\n  [\"945\u00d76=5670\", \"783\u00d75=3915\"],\n  [\"430\u00d78=3440\", \"530\u00d73=1590\"],\n  [\"916\u00d78=7328\", \"660\u00d74=2640\"],\n  [\"808\u00d79=7272\", \"682\u00d78=5456\"],\n  [\"463\u00d79=4167\", \"843\u00d76=5058\"],\n  [\"872\u00d74=3488\", \"972\u00d72=1944\"],\n  [\"730\u00d73=2190\", \"336\u00d78=2688\"],\n  [\"812\u00d73=2436\", \"771\u00d76=4626\"],\n  [\"746\u00d73=2238\", \"511\u00d73=1533\"],\n  [\"552\u00d75=2760\", \"368\u00d79=3312\"],\n  [\"553\u00d73=1659\", \"942\u00d72=1884\"],\n  [\"886\u00d76=5316\", \"740\u00d75=3700\"],\n  [\"570\u00d79=5130\", \"808\u00d73=2424\"],\n  [\"691\u00d73=2073\", \"124\u00d76=744\"],\n  [\"406\u00d78=3248\", \"788\u00d76=4728\"],\n  [\"151\u00d73=453\", \"398\u00d78=3184\"],\n  [\"163\u00d76=978\", \"996\u00d75=4980\"],\n  [\"408\u00d76=2448\", \"215\u00d76=1290\"],\n  [\"722\u00d79=6498\", \"263\u00d75=1315\"],\n  [\"230\u00d75=1150\", \"154\u00d75=770\"],\n  [\"246\u00d74=984\", \"307\u00d77=2149\"],\n  [\"602\u00d76=3612\", \"738\u00d76=4428\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and every three-digit-by-one-digit multiplication\n# answer in the practice table, per the commit's regenerated data.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"),\n    @(\"500\u00d75=2500\", \"338\u00d73=1014\"),\n    @(\"586\u00d72=1172\", \"807\u00d76=4842\"),\n    @(\"341\u00d74=1364\", \"574\u00d77=4018\"),\n    @(\"945\u00d76=5670\", \"783\u00d75=3915\"),\n    @(\"430\u00d78=3440\", \"530\u00d73=1590\"),\n    @(\"916\u00d78=7328\", \"660\u00d74=2640\"),\n    @(\"808\u00d79=7272\", \"682\u00d78=5456\"),\n    @(\"463\u00d79=4167\", \"843\u00d76=5058\"),\n    @(\"872\u00d74=3488\", \"972\u00d72=1944\"),\n    @(\"730\u00d73=2190\", \"336\u00d78=2688\"),\n    @(\"812\u00d73=2436\", \"771\u00d76=4626\"),\n    @(\"746\u00d73=2238\", \"511\u00d73=1533\"),\n    @(\"552\u00d75=2760\", \"368\u00d79=3312\"),\n    @(\"553\u00d73=1659\", \"942\u00d72=1884\"),\n    @(\"886\u00d76=5316\", \"740\u00d75=3700\"),\n    @(\"570\u00d79=5130\", \"808\u00d73=2424\"),\n    @(\"691\u00d73=2073\", \"124\u00d76=744\"),\n    @(\"406\u00d78=3248\", \"788\u00d76=4728\"),\n    @(\"151\u00d73=453\", \"398\u00d78=3184\"),\n    @(\"163\u00d76=978\", \"996\u00d75=4980\"),\n    @(\"408\u00d76=2448\", \"215\u00d76=1290\"),\n    @(\"722\u00d79=6498\", \"263\u00d75=1315\"),\n    @(\"230\u00d75=1150\", \"154\u00d75=770\"),\n    @(\"246\u00d74=984\", \"307\u00d77=2149\"),\n    @(\"602\u00d76=3612\", \"738\u00d76=4428\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
